# ITO-000 Create endpoint user-003
# Adds a new "#user-003" (Get user info by email) endpoint row to the
# "endpoint" sheet, together with a new "Status" / response-message block
# of columns (Status, Message code, Message text) describing the possible
# HTTP responses for every endpoint already documented in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("endpoint")

# 1. Make room for the new "Status" column: insert a column before the
#    existing "Response - DataModel" column (K). This pushes the two
#    existing response columns (Response - DataModel / Response-ResponseBody)
#    one slot to the right (K,L -> L,M).
$ws.Columns.Item(11).Insert()

# 2. Make room for a second response row ("Server error"/500) under the
#    #user-001 endpoint: insert a row before row 7 (pushes the #user-002
#    endpoint row from row 7 down to row 8).
$ws.Rows.Item(7).Insert()

# --- #user-003 (row 9, new endpoint) : basic endpoint fields ---------------
$ws.Cells.Item(9, 1).Value = "#user-003"
$ws.Cells.Item(9, 2).Value = "Get user info by email"
$ws.Cells.Item(9, 3).Value = "/users"
$ws.Cells.Item(9, 4).Value = "Get"
$ws.Cells.Item(9, 5).Value = "email"
$ws.Cells.Item(9, 9).Value = "email"
$ws.Cells.Item(9, 10).Value = "JSON"
$ws.Cells.Item(9, 11).Value = 200
$ws.Cells.Item(9, 12).Value = "onlinejudge.domain.User"

# --- Header row (row 5): new response-message columns N, O ----------------
$ws.Cells.Item(5, 14).Value = "Description response"
$ws.Cells.Item(5, 15).Value = "Message code"
$ws.Cells.Item(5, 16).Value = "Message text"
$ws.Cells.Item(5, 14).WrapText = $true
$ws.Rows.Item(5).RowHeight = 30

# --- #user-003 (row 9): response description --------------------------------
$ws.Cells.Item(9, 14).Value = "find success"

# --- Header row (row 5): new "Status" column K ------------------------------
$ws.Cells.Item(5, 11).Value = "Status"

# --- #user-003 (row 10, new) extra response ---------------------------------
$ws.Cells.Item(10, 11).Value = 400
$ws.Cells.Item(10, 12).Value = "onlinejudge.dto.MyResponse"
$ws.Cells.Item(10, 14).Value = "email not exist"
$ws.Cells.Item(10, 15).Value = "user.email.not.exist"
$ws.Cells.Item(10, 16).Value = "Email [{0}] dose not exist."

# --- #user-001 (row 6) response details --------------------------------------
$ws.Cells.Item(6, 11).Value = 200
$ws.Cells.Item(6, 14).Value = "Create success"

# --- #user-001 (row 7, new) extra response -----------------------------------
$ws.Cells.Item(7, 11).Value = 500
$ws.Cells.Item(7, 14).Value = "Server error"

# --- Column widths for the newly introduced / widened columns --------------
$ws.Columns.Item(11).ColumnWidth = 7.3
$ws.Columns.Item(12).ColumnWidth = 27.85546875
$ws.Columns.Item(14).ColumnWidth = 20.28515625
$ws.Columns.Item(15).ColumnWidth = 23
$ws.Columns.Item(16).ColumnWidth = 23.7109375

# --- Selection mirrors the author's final cursor position ------------------
$ws.Range("O10").Select()
